$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 10-16 already exist; their category data is being reordered/replaced because
# three new "Spiral" sampling schemes were inserted ahead of the rotation/grid schemes,
# and the Gaussian-Quadrature row moved up right after the "Ring Perpendicular" rows.
# Rows 17-19 are brand new rows needed for the schemes that got pushed further down.

$ws.Cells.Item(10,2).Value = "Gaussian-Quadrature"
$ws.Cells.Item(10,3).Value = 1.187930224880939
$ws.Cells.Item(10,4).Value = 0.9440458729845717
$ws.Cells.Item(10,5).Value = 0.9527852964205011
$ws.Cells.Item(10,6).Value = 0.9626495799918713
$ws.Cells.Item(10,7).Value = 1.187930224880939
$ws.Cells.Item(10,8).Value = 0.9440458729845717
$ws.Cells.Item(10,9).Value = 1.022882622533515
$ws.Cells.Item(10,10).Value = 0.9260464903722608
$ws.Cells.Item(10,11).Value = 1.045717215343752
$ws.Cells.Item(10,12).Value = 0.9352641577527503
$ws.Cells.Item(10,13).Value = 1.187930224880939
$ws.Cells.Item(10,14).Value = 0.9484155847025364
$ws.Cells.Item(10,15).Value = 1.011852743569471
$ws.Cells.Item(10,16).Value = 0.9971651825350201

$ws.Cells.Item(11,2).Value = "Spiral-90deg-10rot-5space"
$ws.Cells.Item(11,3).Value = 1.082839528711715
$ws.Cells.Item(11,4).Value = 0.7223566578785771
$ws.Cells.Item(11,5).Value = 1.045892846587956
$ws.Cells.Item(11,6).Value = 0.978169267224371
$ws.Cells.Item(11,7).Value = 1.082839528711715
$ws.Cells.Item(11,8).Value = 0.7223566578785771
$ws.Cells.Item(11,9).Value = 1.057668896531295
$ws.Cells.Item(11,10).Value = 0.9793841194915328
$ws.Cells.Item(11,11).Value = 1.049051808299739
$ws.Cells.Item(11,12).Value = 0.8631387866464122
$ws.Cells.Item(11,13).Value = 1.082839528711715
$ws.Cells.Item(11,14).Value = 0.8841247522332663
$ws.Cells.Item(11,15).Value = 0.9573145751006547
$ws.Cells.Item(11,16).Value = 0.9723127389214498

$ws.Cells.Item(12,2).Value = "Spiral-90deg-15rot-5space"
$ws.Cells.Item(12,3).Value = 1.078273797794782
$ws.Cells.Item(12,4).Value = 0.7241753474120796
$ws.Cells.Item(12,5).Value = 1.046501875024495
$ws.Cells.Item(12,6).Value = 0.9791211063524808
$ws.Cells.Item(12,7).Value = 1.078273797794782
$ws.Cells.Item(12,8).Value = 0.7241753474120796
$ws.Cells.Item(12,9).Value = 1.05682650993508
$ws.Cells.Item(12,10).Value = 0.9808081477001439
$ws.Cells.Item(12,11).Value = 1.047895890025414
$ws.Cells.Item(12,12).Value = 0.8649286022781062
$ws.Cells.Item(12,13).Value = 1.078273797794782
$ws.Cells.Item(12,14).Value = 0.8853386112182874
$ws.Cells.Item(12,15).Value = 0.9570180316459593
$ws.Cells.Item(12,16).Value = 0.9723164095653227

$ws.Cells.Item(13,2).Value = "Spiral-90deg-10rot-3space"
$ws.Cells.Item(13,3).Value = 1.082122467249956
$ws.Cells.Item(13,4).Value = 0.7228950542294115
$ws.Cells.Item(13,5).Value = 1.045572387391628
$ws.Cells.Item(13,6).Value = 0.978495318105336
$ws.Cells.Item(13,7).Value = 1.082122467249956
$ws.Cells.Item(13,8).Value = 0.7228950542294115
$ws.Cells.Item(13,9).Value = 1.057270033813267
$ws.Cells.Item(13,10).Value = 0.9795575103988808
$ws.Cells.Item(13,11).Value = 1.048912861263645
$ws.Cells.Item(13,12).Value = 0.8636386829240843
$ws.Cells.Item(13,13).Value = 1.082122467249956
$ws.Cells.Item(13,14).Value = 0.8842337208105198
$ws.Cells.Item(13,15).Value = 0.9572713067440828
$ws.Cells.Item(13,16).Value = 0.972308039422026

$ws.Cells.Item(14,2).Value = "NoRotation-tilt60deg"
$ws.Cells.Item(14,3).Value = 0.6529599999999995
$ws.Cells.Item(14,4).Value = 0.5621440000000003
$ws.Cells.Item(14,5).Value = 1.571404
$ws.Cells.Item(14,6).Value = 0.8727600000000002
$ws.Cells.Item(14,7).Value = 0.6529599999999995
$ws.Cells.Item(14,8).Value = 0.5621440000000003
$ws.Cells.Item(14,9).Value = 1.247939999999997
$ws.Cells.Item(14,10).Value = 1.141883999999998
$ws.Cells.Item(14,11).Value = 0.8743560000000007
$ws.Cells.Item(14,12).Value = 0.7628520000000003
$ws.Cells.Item(14,13).Value = 0.6529599999999995
$ws.Cells.Item(14,14).Value = 1.066774
$ws.Cells.Item(14,15).Value = 0.914817
$ws.Cells.Item(14,16).Value = 0.9607874999999995

$ws.Cells.Item(15,2).Value = "Rotation-NoTilt"
$ws.Cells.Item(15,3).Value = 0.68
$ws.Cells.Item(15,4).Value = 0.21
$ws.Cells.Item(15,5).Value = 1.979862499999999
$ws.Cells.Item(15,6).Value = 0.696262499999999
$ws.Cells.Item(15,7).Value = 0.68
$ws.Cells.Item(15,8).Value = 0.21
$ws.Cells.Item(15,9).Value = 1.5
$ws.Cells.Item(15,10).Value = 1.15
$ws.Cells.Item(15,11).Value = 0.84
$ws.Cells.Item(15,12).Value = 0.49
$ws.Cells.Item(15,13).Value = 0.68
$ws.Cells.Item(15,14).Value = 1.09493125
$ws.Cells.Item(15,15).Value = 0.8915312499999996
$ws.Cells.Item(15,16).Value = 0.9432656249999998

$ws.Cells.Item(16,2).Value = "Rotation-60detTilt"
$ws.Cells.Item(16,3).Value = 0.8210743283712032
$ws.Cells.Item(16,4).Value = 0.5401972946944007
$ws.Cells.Item(16,5).Value = 1.560636958208003
$ws.Cells.Item(16,6).Value = 0.8248496616447993
$ws.Cells.Item(16,7).Value = 0.8210743283712032
$ws.Cells.Item(16,8).Value = 0.5401972946944007
$ws.Cells.Item(16,9).Value = 1.282742223872001
$ws.Cells.Item(16,10).Value = 1.080514512076795
$ws.Cells.Item(16,11).Value = 0.9089590743040022
$ws.Cells.Item(16,12).Value = 0.7039447246848015
$ws.Cells.Item(16,13).Value = 0.8210743283712032
$ws.Cells.Item(16,14).Value = 1.050417126451202
$ws.Cells.Item(16,15).Value = 0.9366895607296015
$ws.Cells.Item(16,16).Value = 0.9653648472320007

$ws.Cells.Item(17,2).Value = "HexGrid-90degTilt5degRes"
$ws.Cells.Item(17,3).Value = 0.9997038147374451
$ws.Cells.Item(17,4).Value = 0.9980247514947538
$ws.Cells.Item(17,5).Value = 0.992634707969481
$ws.Cells.Item(17,6).Value = 0.9945674475853318
$ws.Cells.Item(17,7).Value = 0.9997038147374451
$ws.Cells.Item(17,8).Value = 0.9980247514947538
$ws.Cells.Item(17,9).Value = 0.9945417256556337
$ws.Cells.Item(17,10).Value = 0.9951981430911249
$ws.Cells.Item(17,11).Value = 0.9960335587309795
$ws.Cells.Item(17,12).Value = 0.994183029362718
$ws.Cells.Item(17,13).Value = 0.9996676384160765
$ws.Cells.Item(17,14).Value = 0.9953297297321173
$ws.Cells.Item(17,15).Value = 0.996232680446753
$ws.Cells.Item(17,16).Value = 0.9956108973284334

$ws.Cells.Item(18,2).Value = "HexGrid-90degTilt22p5degRes"
$ws.Cells.Item(18,3).Value = 0.9357686409085104
$ws.Cells.Item(18,4).Value = 0.9768712914227272
$ws.Cells.Item(18,5).Value = 1.024392138926677
$ws.Cells.Item(18,6).Value = 1.004445348355196
$ws.Cells.Item(18,7).Value = 0.9357686409085104
$ws.Cells.Item(18,8).Value = 0.9768712914227272
$ws.Cells.Item(18,9).Value = 1.000406596344744
$ws.Cells.Item(18,10).Value = 1.019365122437331
$ws.Cells.Item(18,11).Value = 0.984337981661255
$ws.Cells.Item(18,12).Value = 0.9959729712081115
$ws.Cells.Item(18,13).Value = 0.9357686409085104
$ws.Cells.Item(18,14).Value = 1.000631715174702
$ws.Cells.Item(18,15).Value = 0.9853693549032776
$ws.Cells.Item(18,16).Value = 0.9926950114080689

$ws.Cells.Item(19,2).Value = "HexGrid-60degTilt5degRes"
$ws.Cells.Item(19,3).Value = 0.9707473758416045
$ws.Cells.Item(19,4).Value = 1.072239489946939
$ws.Cells.Item(19,5).Value = 0.9850136336673717
$ws.Cells.Item(19,6).Value = 1.002521268804754
$ws.Cells.Item(19,7).Value = 0.9707473758416045
$ws.Cells.Item(19,8).Value = 1.072239489946939
$ws.Cells.Item(19,9).Value = 0.9723512690807741
$ws.Cells.Item(19,10).Value = 0.998872779636711
$ws.Cells.Item(19,11).Value = 0.9788047132915261
$ws.Cells.Item(19,12).Value = 1.036646142126203
$ws.Cells.Item(19,13).Value = 0.9707473758416045
$ws.Cells.Item(19,14).Value = 1.028626561807155
$ws.Cells.Item(19,15).Value = 1.007630442065167
$ws.Cells.Item(19,16).Value = 1.002149584049485

# New rows 17-19 need column A populated with the row index and the same
# bold/centered/bordered style that the rest of column A already uses.
$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(18,1).Value = 16
$ws.Cells.Item(19,1).Value = 17

$ws.Range("A16").Copy()
$ws.Range("A17:A19").PasteSpecial(-4122)
$excel.CutCopyMode = $false
